$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "330.99"
Set-TextValue $ws.Range("E2") "-0.15%"
Set-TextValue $ws.Range("D3") "41.60"
Set-TextValue $ws.Range("E3") "0.53%"
Set-TextValue $ws.Range("D4") "5.673"
Set-TextValue $ws.Range("E4") "-1.38%"
Set-TextValue $ws.Range("D5") "0.08352"
Set-TextValue $ws.Range("E5") "2.91%"
Set-TextValue $ws.Range("D6") "8.792"
Set-TextValue $ws.Range("E6") "0.95%"
Set-TextValue $ws.Range("D7") "2.005"
Set-TextValue $ws.Range("E7") "-5.14%"
Set-TextValue $ws.Range("D8") "4.512"
Set-TextValue $ws.Range("E8") "0.22%"
Set-TextValue $ws.Range("D9") "2.938"
Set-TextValue $ws.Range("E9") "-1.34%"
Set-TextValue $ws.Range("D10") "0.9266"
Set-TextValue $ws.Range("E10") "0.06%"
Set-TextValue $ws.Range("E11") "1.00%"
Set-TextValue $ws.Range("D12") "0.1963"
Set-TextValue $ws.Range("E12") "0.47%"
Set-TextValue $ws.Range("D13") "0.09375"
Set-TextValue $ws.Range("E13") "2.37%"
Set-TextValue $ws.Range("D14") "0.03886"
Set-TextValue $ws.Range("E14") "6.45%"
Set-TextValue $ws.Range("D15") "0.1060"
Set-TextValue $ws.Range("E15") "0.89%"
Set-TextValue $ws.Range("D16") "0.001307"
Set-TextValue $ws.Range("E16") "0.63%"
Set-TextValue $ws.Range("D17") "0.006115"
Set-TextValue $ws.Range("E17") "-3.25%"
Set-TextValue $ws.Range("D18") "3.441"
Set-TextValue $ws.Range("E18") "2.01%"
Set-TextValue $ws.Range("E19") "1.24%"
Set-TextValue $ws.Range("D20") "8.472"
Set-TextValue $ws.Range("E20") "-4.31%"
Set-TextValue $ws.Range("E21") "-0.91%"
Set-TextValue $ws.Range("D22") "0.2488"
Set-TextValue $ws.Range("E22") "-4.42%"
Set-TextValue $ws.Range("D23") "0.04414"
Set-TextValue $ws.Range("E23") "-0.13%"
Set-TextValue $ws.Range("D24") "0.001267"
Set-TextValue $ws.Range("E24") "0.94%"
Set-TextValue $ws.Range("D25") "0.004401"
Set-TextValue $ws.Range("E25") "-2.77%"
Set-TextValue $ws.Range("D26") "0.0001201"
Set-TextValue $ws.Range("E26") "-3.02%"
Set-TextValue $ws.Range("D39") "0.02805"
Set-TextValue $ws.Range("E39") "1.35%"
Set-TextValue $ws.Range("D40") "0.05549"
Set-TextValue $ws.Range("E40") "0.33%"
Set-TextValue $ws.Range("D41") "0.007798"
Set-TextValue $ws.Range("E41") "2.39%"
Set-TextValue $ws.Range("D42") "0.1436"
Set-TextValue $ws.Range("E42") "0.71%"
Set-TextValue $ws.Range("D43") "0.009298"
Set-TextValue $ws.Range("E43") "-5.62%"
Set-TextValue $ws.Range("D44") "0.002101"
Set-TextValue $ws.Range("E44") "-5.38%"
Set-TextValue $ws.Range("D45") "0.01065"
Set-TextValue $ws.Range("E45") "-9.86%"
Set-TextValue $ws.Range("D46") "0.00006992"
Set-TextValue $ws.Range("E46") "3.24%"
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.21%"
Set-TextValue $ws.Range("D48") "0.003546"
Set-TextValue $ws.Range("E48") "15.61%"
Set-TextValue $ws.Range("D49") "0.002280"
Set-TextValue $ws.Range("E49") "0.18%"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "0.21%"
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "0.21%"
